$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the unfinished "kalecaesar" placeholder item (row 3, column F)
# with the "placeholder" text, matching the commit message:
# "Cleared empty cells and replaced all unfinished items with placeholders"
$ws.Range("F3").Value = "placeholder"

# Update the active selection to match the saved workbook state.
$ws.Range("B10").Select()
